$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 text
$ws.Range("B2").Value = "cuak_v1-1.95"
$ws.Range("C2").Value = "The first containerized Chrome deployments (Multi-VM behind ALB, VNC session routing issue)"

# Add new row 3 (45738 = 3/22/2025 as an Excel date serial)
$ws.Range("A3").Value = 45738
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B3").Value = "cuak_v2"
$ws.Range("C3").Value = "FIRST STABLE RELEASE (Implements Redis Caching to find VM running Chrome container)"

# Column C grows to fit the longer note text
$ws.Columns.Item(3).ColumnWidth = 79.8

# Move selection
$ws.Range("A4").Select()
